$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "The cat is on the roof."
$ws.Range("B2").Value = "Die Katze ist auf dem Dach."
$ws.Range("C2").Value = "Katten är på taket."

$ws.Range("A3").Value = "Signal improvement"
$ws.Range("B3").Value = "Signalverbesserung"
$ws.Range("C3").Value = "Signalförbättring"
